# "Fixed mistakes in precedence table" - five cells in the precedence
# table incorrectly showed '>' where they should show '<'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<"
$ws.Range("C2").Value = "<"
$ws.Range("C3").Value = "<"
$ws.Range("G9").Value = "<"
$ws.Range("H9").Value = "<"

# View state at save time: zoomed in further, selection moved to T4.
$excel.ActiveWindow.Zoom = 115
$ws.Range("T4").Select()
